$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (the "sdmx-dimension" row),
# shifting the existing rows 2-5 down to rows 3-6.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the machine-readable column
# identifiers that relate to the human-readable headers in row 1.
$ws.Range("A2").Value = "ccaa-nombre"
$ws.Range("B2").Value = "siglas"
$ws.Range("C2").Value = "ccaa-codigo"
$ws.Range("D2").Value = "diputados"
$ws.Range("E2").Value = "provincia-codigo"
$ws.Range("F2").Value = "provincia-nombre"
$ws.Range("G2").Value = "ano"
$ws.Range("H2").Value = "votos"

# The former row 5 (now shifted to row 6) only contained a stray
# "mapping-ano.xlsx" value in column G and is no longer needed now that
# the previous row 4 (now row 5) carries the complete xsd:date entry.
$ws.Rows.Item(6).Delete()
